$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.1
$ws.Range("G2").Value = 4.2
$ws.Range("H2").Value = 2.5
$ws.Range("I2").Value = 2.8
$ws.Range("L2").Value = 1.57
$ws.Range("Q2").Value = 2.74
$ws.Range("R2").Value = 1.13
$ws.Range("V2").Value = 1.55

# Row 3
$ws.Range("F3").Value = 1.16
$ws.Range("H3").Value = 16.5
$ws.Range("I3").Value = 21
$ws.Range("J3").Value = 9
$ws.Range("L3").Value = 1.16
$ws.Range("R3").Value = 2.04
$ws.Range("S3").Value = 1.78
$ws.Range("T3").Value = 1.96
$ws.Range("W3").Value = 5.8
$ws.Range("X3").Value = 990
$ws.Range("Y3").Value = 990
$ws.Range("AB3").Value = 18.5
$ws.Range("AC3").Value = 28
$ws.Range("AG3").Value = 16.5
$ws.Range("AH3").Value = 44
$ws.Range("AL3").Value = 980
$ws.Range("AN3").Value = 2.8

# Row 4
$ws.Range("F4").Value = 1.83
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 5.1
$ws.Range("I4").Value = 6.2
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 3.6
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 2.56
$ws.Range("O4").Value = 1.52
$ws.Range("P4").Value = 1.52
$ws.Range("Q4").Value = 2.44
$ws.Range("R4").Value = 1.18
$ws.Range("S4").Value = 4.6
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 1.67
$ws.Range("V4").Value = 1.19
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 10.5
$ws.Range("Y4").Value = 14.5
$ws.Range("Z4").Value = 44
$ws.Range("AB4").Value = 6.6
$ws.Range("AC4").Value = 9.4
$ws.Range("AD4").Value = 25
$ws.Range("AE4").Value = 140
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 13.5
$ws.Range("AH4").Value = 28
$ws.Range("AJ4").Value = 24
$ws.Range("AK4").Value = 28
$ws.Range("AL4").Value = 75
$ws.Range("AN4").Value = 24

# Row 5
$ws.Range("F5").Value = 2.56
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 3.4
$ws.Range("X5").Value = 10
$ws.Range("AB5").Value = 8.800000000000001

# Row 6
$ws.Range("F6").Value = 1.65
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 4.8
$ws.Range("N6").Value = 7.8
$ws.Range("O6").Value = 1.11
$ws.Range("P6").Value = 3.3
$ws.Range("R6").Value = 1.94
$ws.Range("S6").Value = 1.87
$ws.Range("V6").Value = 1.25
$ws.Range("Y6").Value = 980
$ws.Range("AE6").Value = 48
$ws.Range("AL6").Value = 23
$ws.Range("AO6").Value = 28

# Row 7
$ws.Range("G7").Value = 5.6
$ws.Range("N7").Value = 3.5
$ws.Range("P7").Value = 1.87
$ws.Range("Q7").Value = 1.93
$ws.Range("T7").Value = 1.83
$ws.Range("U7").Value = 1.97
$ws.Range("Y7").Value = 980
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 980
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 980
$ws.Range("AF7").Value = 980
$ws.Range("AG7").Value = 980
$ws.Range("AH7").Value = 980
$ws.Range("AI7").Value = 980
$ws.Range("AO7").Value = 1000

# Row 8
$ws.Range("F8").Value = 1.62
$ws.Range("G8").Value = 1.76
$ws.Range("H8").Value = 6.6
$ws.Range("J8").Value = 3.2
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 2.7
$ws.Range("O8").Value = 1.47
$ws.Range("Q8").Value = 2.38
$ws.Range("S8").Value = 4.8
$ws.Range("W8").Value = 2.3

# Row 9
$ws.Range("V9").Value = 1.1

# Row 10
$ws.Range("F10").Value = 1.63
$ws.Range("G10").Value = 1.64
$ws.Range("R10").Value = 1.45
$ws.Range("W10").Value = 2.56
$ws.Range("X10").Value = 18
$ws.Range("AA10").Value = 170
$ws.Range("AG10").Value = 9.4
$ws.Range("AN10").Value = 8.4
$ws.Range("AO10").Value = 85

# Row 11
$ws.Range("G11").Value = 2.52
$ws.Range("I11").Value = 5
$ws.Range("V11").Value = 1.25
$ws.Range("W11").Value = 1.66
